$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Part 1: merge the two runs "THU Sep 28" + " 14:05:43 PDT 2017" into
# a single run "THU Sep 28 14:05:43 PDT 2017" (same visible text,
# Find/Replace naturally collapses the match into one run).
# -------------------------------------------------------------------
$d.Content.Find.Execute("THU Sep 28 14:05:43 PDT 2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "THU Sep 28 14:05:43 PDT 2017", 2) | Out-Null

# -------------------------------------------------------------------
# Part 2: append a new "purchase details" entry (SAT SEP 30) after the
# very last existing entry (the one ending "...Amount balance ... -
# 284312.0").
# -------------------------------------------------------------------

# Locate the anchor paragraph: the last paragraph whose text contains
# "284312.0" (the final "Amount balance" line in the document).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*- 284312.0*") {
        $anchorIndex = $i
    }
}

$cur = $d.Paragraphs.Item($anchorIndex)

# --- New paragraph 1: empty, bold (separator line) ---------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 1)
$cur.Range.Delete()

# --- New paragraph 2: date line, not bold -------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 2)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("SAT SEP 30 16:38:41 PDT 2017")

# --- New paragraph 3: Person Name - NG ----------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 3)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("Person Name`t`t`t`t- NG")

# --- New paragraph 4: dashed separator -----------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 4)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("---------------------------------------------------------------")

# --- New paragraph 5: Item Name - CARROT 1 -------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 5)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("Item Name`t`t`t`t- CARROT 1")

# --- New paragraph 6: Number of Pockets - 8 ------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 6)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("Number of Pockets`t`t`t- 8")

# --- New paragraph 7: Number of KGs - 775 --------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 7)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("Number of KGs`t`t`t- 775")

# --- New paragraph 8: Rate - 20 ------------------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 8)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("Rate`t`t`t`t`t- 20")

# --- New paragraph 9: Transport & Miscellaneous - 1320 -------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 9)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("Transport & Miscellaneous`t- 1320")

# --- New paragraph 10: Total Price - 16820.0 ------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 10)
$cur.Range.Font.Bold = 0
$cur.Range.InsertAfter("Total Price`t`t`t`t- 16820.0")

# --- New paragraph 11: Amount balance - 301132.0 (bold) -------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 11)
$cur.Range.Font.Bold = 1
$cur.Range.InsertAfter("Amount balance`t`t`t- 301132.0")

# --- New paragraph 12: empty, not bold -------------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 12)
$cur.Range.Font.Bold = 0
$cur.Range.Delete()

# --- New paragraph 13: empty, bold -----------------------------------------
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Item($anchorIndex + 13)
$cur.Range.Delete()

Write-Output "done"
